$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.317.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.182.29"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.56%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.70%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -5.34%  "

$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.180.99"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.31%  "

$ws.Range("E10").Value = "  -3.72%  "

$ws.Range("E11").Value = "  -0.92%  "

$ws.Range("E12").Value = "  -5.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.729.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.39%  "

$ws.Range("E14").Value = "  -1.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.389.25"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000158"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.96%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.180.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "419.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("E22").Value = "  -4.38%  "

$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.41%  "

$ws.Range("E25").Value = "  +2.78%  "

$ws.Range("E26").Value = "  -4.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000107"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.85%  "

$ws.Range("E28").Value = "  -2.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("E32").Value = "  -0.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.02%  "

$ws.Range("E35").Value = "  -4.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "155.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.07%  "

$ws.Range("E37").Value = "  -5.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.704.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.65%  "

$ws.Range("E39").Value = "  -6.73%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.75%  "

$ws.Range("E41").Value = "  -3.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.27%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.709"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0625"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.72%  "

$ws.Range("E46").Value = "  -5.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "297.07"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.76%  "

$ws.Range("E48").Value = "  -2.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.04"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.67%  "

$ws.Range("E50").Value = "  -4.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.02%  "
